$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @("2023-12-07 20:28:00", 0.001),
    @("2023-12-07 20:29:21", 0.0044),
    @("2023-12-07 20:30:11", 0.003),
    @("2023-12-07 20:30:18", 0.0004)
)

$startRow = 85
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
